# Helper table for poke correlation calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing per-rat infusion counts in column E (row 7-14) ---
$ws.Range("E7").Value  = 38
$ws.Range("E8").Value  = 22
$ws.Range("E9").Value  = 29
$ws.Range("E10").Value = 31
$ws.Range("E11").Value = 22
$ws.Range("E12").Value = 16
$ws.Range("E13").Value = 4
$ws.Range("E14").Value = 28

# --- Move the footer note from row 17 down to row 20 ---
$ws.Range("E17").ClearContents()
$ws.Range("E20").Value = "Data are number of infusions/6 hr session"

# --- Add new helper rows 15-17 (Ken, Woody, Trixie) ---
$ws.Range("D15").Value = "Ken"
$ws.Range("E15").Value = 20
$ws.Range("F15").Value = 14
$ws.Range("G15").Value = 10
$ws.Range("H15").Value = 36
$ws.Range("I15").Value = 40
$ws.Range("J15").Value = 66
$ws.Range("K15").Value = 32
$ws.Range("L15").Value = 107
$ws.Range("M15").Value = 77
$ws.Range("N15").Value = "nan"
$ws.Range("O15").Value = "nan"
$ws.Range("P15").Value = "nan"
$ws.Range("Q15").Value = 21
$ws.Range("R15").Value = 61

$ws.Range("D16").Value = "Woody"
$ws.Range("E16").Value = 31
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 9
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = 14
$ws.Range("L16").Value = 14
$ws.Range("M16").Value = 10
$ws.Range("N16").Value = 9
$ws.Range("O16").Value = 10
$ws.Range("P16").Value = 23.5
$ws.Range("Q16").Value = 37
$ws.Range("R16").Value = 37

$ws.Range("D17").Value = "Trixie"
$ws.Range("E17").Value = 21
$ws.Range("F17").Value = 17.5
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 24
$ws.Range("I17").Value = 23
$ws.Range("J17").Value = 31
$ws.Range("K17").Value = 7
$ws.Range("L17").Value = 10
$ws.Range("M17").Value = 8.5
$ws.Range("N17").Value = 7
$ws.Range("O17").Value = 19
$ws.Range("P17").Value = 17
$ws.Range("Q17").Value = 24
$ws.Range("R17").Value = 3

# --- Update the selection to match the recorded cursor position ---
[void]$ws.Range("K18").Select()
